# The workbook has a raw data table in A1:E43 on Sheet1. Column D is
# currently unused (blank) and the "tool life" column lives in E. The
# edit shifts that "tool life" column left into D (closing the gap),
# fixes two transposed Axial-depth values, widens column H for a note,
# and leaves the selection where the author left it before saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "tool life" column (header + all data rows) from E to D so
# the used range becomes A1:D43 instead of A1:E43.
$ws.Range("E1:E43").Cut($ws.Range("D1:D43"))

# Correct two pairs of swapped "Axial depth" readings (rows 29/30 and
# rows 32/33 each had their values transposed).
$ws.Range("C29").Value = 1
$ws.Range("C30").Value = 1.5
$ws.Range("C32").Value = 1
$ws.Range("C33").Value = 1.5

# Widen column H to fit an added note/comment.
$ws.Range("H1").ColumnWidth = 33.141183035714285

# Restore the selection/viewport the author had active when the file
# was last saved.
$excel.Goto($ws.Range("E1:K43"))
